$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 30 and 31 had their match data swapped (home/away + odds +
#    timestamps + url). Columns A:E (index, country, league, season,
#    kickoff date) are identical between the two rows, so only F:V
#    need to be rewritten to their corrected values.
# ------------------------------------------------------------------
# --- row 30 (Genk 0 - 0 Charleroi) ---
$ws.Range("F30").Value = "Genk"
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = "Charleroi"
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1.47
$ws.Range("K30").Value = "13/08/2023 18:42"
$ws.Range("L30").Value = 1.59
$ws.Range("M30").Value = "20/08/2023 15:53"
$ws.Range("N30").Value = 4.85
$ws.Range("O30").Value = "13/08/2023 18:42"
$ws.Range("P30").Value = 4.57
$ws.Range("Q30").Value = "20/08/2023 15:53"
$ws.Range("R30").Value = 5.52
$ws.Range("S30").Value = "13/08/2023 18:42"
$ws.Range("T30").Value = 5.25
$ws.Range("U30").Value = "20/08/2023 15:58"
$ws.Range("V30").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-charleroi/xlfRPeMc/"

# --- row 31 (Gent 2 - 2 St. Truiden) ---
$ws.Range("F31").Value = "Gent"
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = "St. Truiden"
$ws.Range("I31").Value = 2
$ws.Range("J31").Value = 1.4
$ws.Range("K31").Value = "13/08/2023 19:42"
$ws.Range("L31").Value = 1.52
$ws.Range("M31").Value = "20/08/2023 15:57"
$ws.Range("N31").Value = 5
$ws.Range("O31").Value = "13/08/2023 19:42"
$ws.Range("P31").Value = 4.45
$ws.Range("Q31").Value = "20/08/2023 15:59"
$ws.Range("R31").Value = 6.31
$ws.Range("S31").Value = "13/08/2023 19:42"
$ws.Range("T31").Value = 6.41
$ws.Range("U31").Value = "20/08/2023 15:59"
$ws.Range("V31").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-st-truiden/SxmvNg6G/"

# ------------------------------------------------------------------
# 2) Seven new match rows (56-62) are appended at the bottom of the
#    sheet. Copy formatting (bold/border style on col A, datetime
#    number format on col E) from an existing data row first, then
#    fill in the values.
# ------------------------------------------------------------------
$ws.Range("A2:V2").Copy()
$ws.Range("A56:V56").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A57:V57").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A58:V58").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A59:V59").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A60:V60").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A61:V61").PasteSpecial(-4122)
$ws.Range("A2:V2").Copy()
$ws.Range("A62:V62").PasteSpecial(-4122)

$ws.Range("A56").Value = 55
$ws.Range("B56").Value = "belgium"
$ws.Range("C56").Value = "jupiler-pro-league"
$ws.Range("D56").Value = "2023-2024"
$ws.Range("E56").Value = 45192.66666666666
$ws.Range("F56").Value = "Charleroi"
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = "Kortrijk"
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1.68
$ws.Range("K56").Value = "17/09/2023 15:12"
$ws.Range("L56").Value = 1.51
$ws.Range("M56").Value = "23/09/2023 15:53"
$ws.Range("N56").Value = 4.17
$ws.Range("O56").Value = "17/09/2023 15:12"
$ws.Range("P56").Value = 4.79
$ws.Range("Q56").Value = "23/09/2023 15:58"
$ws.Range("R56").Value = 4.71
$ws.Range("S56").Value = "17/09/2023 15:12"
$ws.Range("T56").Value = 6.08
$ws.Range("U56").Value = "23/09/2023 15:58"
$ws.Range("V56").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/charleroi-kortrijk/GxrGORbI/"

$ws.Range("A57").Value = 56
$ws.Range("B57").Value = "belgium"
$ws.Range("C57").Value = "jupiler-pro-league"
$ws.Range("D57").Value = "2023-2024"
$ws.Range("E57").Value = 45192.76041666666
$ws.Range("F57").Value = "KV Mechelen"
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = "Leuven"
$ws.Range("I57").Value = 2
$ws.Range("J57").Value = 2
$ws.Range("K57").Value = "17/09/2023 18:43"
$ws.Range("L57").Value = 1.91
$ws.Range("M57").Value = "23/09/2023 18:14"
$ws.Range("N57").Value = 3.57
$ws.Range("O57").Value = "17/09/2023 18:43"
$ws.Range("P57").Value = 4.01
$ws.Range("Q57").Value = "23/09/2023 18:14"
$ws.Range("R57").Value = 3.55
$ws.Range("S57").Value = "17/09/2023 18:43"
$ws.Range("T57").Value = 3.83
$ws.Range("U57").Value = "23/09/2023 18:13"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/kv-mechelen-leuven/rPq8QmU5/"

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = "belgium"
$ws.Range("C58").Value = "jupiler-pro-league"
$ws.Range("D58").Value = "2023-2024"
$ws.Range("E58").Value = 45192.86458333334
$ws.Range("F58").Value = "Antwerp"
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = "RWDM"
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1.36
$ws.Range("K58").Value = "16/09/2023 17:13"
$ws.Range("L58").Value = 1.41
$ws.Range("M58").Value = "23/09/2023 20:40"
$ws.Range("N58").Value = 5.1
$ws.Range("O58").Value = "16/09/2023 17:13"
$ws.Range("P58").Value = 5.14
$ws.Range("Q58").Value = "23/09/2023 20:44"
$ws.Range("R58").Value = 8.73
$ws.Range("S58").Value = "16/09/2023 17:13"
$ws.Range("T58").Value = 7.54
$ws.Range("U58").Value = "23/09/2023 20:44"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/antwerp-rwd-molenbeek/UoBaUVUt/"

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "belgium"
$ws.Range("C59").Value = "jupiler-pro-league"
$ws.Range("D59").Value = "2023-2024"
$ws.Range("E59").Value = 45193.5625
$ws.Range("F59").Value = "Genk"
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = "St. Truiden"
$ws.Range("I59").Value = 3
$ws.Range("J59").Value = 1.53
$ws.Range("K59").Value = "17/09/2023 18:43"
$ws.Range("L59").Value = 1.56
$ws.Range("M59").Value = "24/09/2023 13:04"
$ws.Range("N59").Value = 4.39
$ws.Range("O59").Value = "17/09/2023 18:43"
$ws.Range("P59").Value = 4.55
$ws.Range("Q59").Value = "24/09/2023 13:05"
$ws.Range("R59").Value = 6
$ws.Range("S59").Value = "17/09/2023 18:43"
$ws.Range("T59").Value = 5.71
$ws.Range("U59").Value = "24/09/2023 13:05"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-st-truiden/xCE7S9ag/"

$ws.Range("A60").Value = 59
$ws.Range("B60").Value = "belgium"
$ws.Range("C60").Value = "jupiler-pro-league"
$ws.Range("D60").Value = "2023-2024"
$ws.Range("E60").Value = 45193.66666666666
$ws.Range("F60").Value = "Cercle Brugge KSV"
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = "Royale Union SG"
$ws.Range("I60").Value = 2
$ws.Range("J60").Value = 2.55
$ws.Range("K60").Value = "16/09/2023 17:43"
$ws.Range("L60").Value = 2.57
$ws.Range("M60").Value = "24/09/2023 15:59"
$ws.Range("N60").Value = 3.4
$ws.Range("O60").Value = "16/09/2023 17:43"
$ws.Range("P60").Value = 3.66
$ws.Range("Q60").Value = "24/09/2023 15:59"
$ws.Range("R60").Value = 2.67
$ws.Range("S60").Value = "16/09/2023 17:43"
$ws.Range("T60").Value = 2.7
$ws.Range("U60").Value = "24/09/2023 15:59"
$ws.Range("V60").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/cercle-brugge-royale-union-sg/rkF3Tkpm/"

$ws.Range("A61").Value = 60
$ws.Range("B61").Value = "belgium"
$ws.Range("C61").Value = "jupiler-pro-league"
$ws.Range("D61").Value = "2023-2024"
$ws.Range("E61").Value = 45193.77083333334
$ws.Range("F61").Value = "Anderlecht"
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = "Club Brugge KV"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 2.58
$ws.Range("K61").Value = "17/09/2023 15:12"
$ws.Range("L61").Value = 3.82
$ws.Range("M61").Value = "24/09/2023 18:29"
$ws.Range("N61").Value = 3.44
$ws.Range("O61").Value = "17/09/2023 15:12"
$ws.Range("P61").Value = 4.01
$ws.Range("Q61").Value = "24/09/2023 18:29"
$ws.Range("R61").Value = 2.62
$ws.Range("S61").Value = "17/09/2023 15:12"
$ws.Range("T61").Value = 1.91
$ws.Range("U61").Value = "24/09/2023 18:29"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/anderlecht-club-brugge/G2DBRTEa/"

$ws.Range("A62").Value = 61
$ws.Range("B62").Value = "belgium"
$ws.Range("C62").Value = "jupiler-pro-league"
$ws.Range("D62").Value = "2023-2024"
$ws.Range("E62").Value = 45193.80208333334
$ws.Range("F62").Value = "Gent"
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = "Eupen"
$ws.Range("I62").Value = 1
$ws.Range("J62").Value = 1.29
$ws.Range("K62").Value = "17/09/2023 17:42"
$ws.Range("L62").Value = 1.23
$ws.Range("M62").Value = "24/09/2023 19:11"
$ws.Range("N62").Value = 5.45
$ws.Range("O62").Value = "17/09/2023 17:42"
$ws.Range("P62").Value = 6.81
$ws.Range("Q62").Value = "24/09/2023 19:11"
$ws.Range("R62").Value = 8.83
$ws.Range("S62").Value = "17/09/2023 17:42"
$ws.Range("T62").Value = 11.7
$ws.Range("U62").Value = "24/09/2023 19:11"
$ws.Range("V62").Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-eupen/llcVob0E/"

